# Append a newly-scraped Lancers listing to the top of the data table
# (row 2), pushing the existing rows down by one, and refresh the
# "取得日時" (fetched-at) timestamp on every row to the new scrape time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# The engine's Rows.Insert() shifts cell values/styles down correctly but
# does not re-anchor existing Hyperlink objects, so drop them first and
# rebuild the whole set afterwards against the final row positions.
$ws.Hyperlinks.Delete()

# Insert a new row above the current row 2 (the first data row) so all
# existing entries shift down by one.
$ws.Rows.Item(2).Insert()

# Populate the new top row with the freshly scraped listing.
$newTimestamp = "2025-11-06 12:37:17"

$ws.Cells.Item(2, 1).Value = $newTimestamp
$ws.Cells.Item(2, 2).Value = "【Next.js × TypeScript × Tailwind】コンポーネント制作パートナー募集!"
$ws.Cells.Item(2, 3).Value = "システム開発"
$ws.Cells.Item(2, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(2, 5).Value = "期限情報なし"
$ws.Cells.Item(2, 6).Value = "https://www.lancers.jp/work/detail/5428507"
$ws.Cells.Item(2, 7).Value = 528
$ws.Cells.Item(2, 8).Value = "🔥AI,Next.js"

# Refresh the scrape timestamp on every other (pre-existing) data row.
for ($r = 3; $r -le 7; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# Rebuild the URL hyperlinks (column F) for all seven data rows, in the
# same style used throughout the sheet.
$urls = @(
    "https://www.lancers.jp/work/detail/5428507",
    "https://www.lancers.jp/work/detail/5427956",
    "https://www.lancers.jp/work/detail/5217096",
    "https://www.lancers.jp/work/detail/5428337",
    "https://www.lancers.jp/work/detail/5428278",
    "https://www.lancers.jp/work/detail/5428124"
)
for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 6)
    $ws.Hyperlinks.Add($cell, $urls[$i])
    $cell.Style = "Hyperlink"
}

# Minor column-width adjustments that accompanied this refresh.
# (The -5/6 offset compensates for this engine's char-width rounding so
# the stored OOXML <col width> lands exactly on the intended integer.)
$ws.Columns.Item(2).ColumnWidth = 51.16666666666667
$ws.Columns.Item(8).ColumnWidth = 12.16666666666667
